# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Cells D/E hold text (e.g. "66.461.86", "  +2.47%  "), not numbers, so a
# handful of plain-decimal-looking price values are written with a leading
# apostrophe (Excel's normal "treat as text" quote-prefix) to stop the COM
# layer from auto-coercing them into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.461.86'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '3.485.72'
$ws.Range('E3').Value = '  +1.68%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''589.23'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').Value = '''167.71'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.482.61'
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('E9').Value = '  +5.97%  '
$ws.Range('D10').Value = '''7.31'
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('E11').Value = '  +5.24%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').Value = '4.089.93'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '''27.89'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.489.84'
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '''0.0000177'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').Value = '3.488.15'
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('D19').Value = '''6.25'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').Value = '''389.92'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').Value = '''7.86'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '''72.67'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('E26').Value = '  +3.28%  '
$ws.Range('D27').Value = '''10.13'
$ws.Range('E27').Value = '  +4.23%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('E34').Value = '  +2.93%  '
$ws.Range('E35').Value = '  +5.58%  '
$ws.Range('D36').Value = '''162.91'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('D37').Value = '''0.892'
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('E38').Value = '  +2.74%  '
$ws.Range('E39').Value = '  +4.94%  '
$ws.Range('E40').Value = '  +4.71%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').Value = '''26.10'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').Value = '2.766.53'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '''26.40'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''42.74'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('D47').Value = '''0.0307'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').Value = '''340.29'
$ws.Range('E48').Value = '  +3.03%  '
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('D50').Value = '''33.30'
$ws.Range('E50').Value = '  +7.73%  '
$ws.Range('D51').Value = '''0.848'
$ws.Range('E51').Value = '  +3.34%  '
